# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.461.15"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "2.676.77"
$ws.Range("E3").Value = "  +4.12%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.10"
$ws.Range("E5").Value = "  +5.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.09"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "2.677.51"
$ws.Range("E9").Value = "  +4.14%  "
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.362"
$ws.Range("E13").Value = "  +3.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.41"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "3.158.83"
$ws.Range("E15").Value = "  +4.06%  "
$ws.Range("D16").Value = "63.311.77"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000146"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "2.679.35"
$ws.Range("E18").Value = "  +4.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.47"
$ws.Range("E19").Value = "  +3.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.78"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.42"
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.88"
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.48"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  +3.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.55"
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.71"
$ws.Range("E27").Value = "  +5.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.165"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "544.24"
$ws.Range("E29").Value = "  +18.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.89"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  +7.47%  "
$ws.Range("E33").Value = "  +8.42%  "
$ws.Range("D34").Value = "0.0₃0810"
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "172.93"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.17"
$ws.Range("E36").Value = "  +14.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.406"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  +2.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.86"
$ws.Range("E40").Value = "  +10.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.42"
$ws.Range("E41").Value = "  +11.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.76"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.28"
$ws.Range("E44").Value = "  +5.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0573"
$ws.Range("E45").Value = "  +7.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.637"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0241"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0965"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("E49").Value = "  +5.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.76"
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.32"
$ws.Range("E51").Value = "  -0.73%  "
